$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (pushes "Versandpauschale" row and the
# totals row down by one), mirroring the Excel UI action of inserting a
# table row above the last data row.
$ws.Rows(20).Insert()

# Fix up the part-number counter of the row that got pushed down
# (item 14 "Versandpauschale" becomes item 15).
$ws.Range("B21").Value = 15

# Populate the newly inserted row 20 with the new BOM line item.
$ws.Range("B20").Value = 14
$ws.Range("C20").Value = "Stromkabel"
$ws.Range("D20").Value = "USB zu DC 5,5 * 2,5mm Stromkabel"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.83

# Add the hyperlink for the new part and restore the normal "Link" look
# (adding a hyperlink stamps a fresh style on the cell, so re-apply the
# existing Link cell style afterwards to match the other rows).
$ws.Hyperlinks.Add($ws.Range("H20"), "https://de.aliexpress.com/item/4000993209004.html") | Out-Null
$ws.Range("H20").Value = "https://de.aliexpress.com/item/4000993209004.html"
$ws.Range("H20").Style = "Link"

# Grow the table (Tabelle1) so the new row participates in the table,
# autofilter and the SUM(Tabelle1[Preis]) total.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B6:H21"))

# Recalculate so the total in the (now shifted) summary row picks up the
# new price.
$excel.Calculate()

# Match the saved selection from the source edit.
$ws.Range("C24").Select()
